$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 461.68182
$ws.Range("I107").Value = 444.23077
$ws.Range("J107").Value = 486.8889
$ws.Range("K107").Value = 444.23077
$ws.Range("L107").Value = 486.8889
$ws.Range("M107").Value = 1475.76923
$ws.Range("N107").Value = -4326.8889
# Row 112
$ws.Range("H112").Value = 1061.7142
$ws.Range("I112").Value = 550
$ws.Range("J112").Value = 1092.7273
$ws.Range("K112").Value = 1650
$ws.Range("L112").Value = 3278.1819
$ws.Range("M112").Value = -542
$ws.Range("N112").Value = -5494.1819
# Row 113
$ws.Range("H113").Value = 45458348
$ws.Range("I113").Value = 83336580
$ws.Range("J113").Value = 4478.7
$ws.Range("K113").Value = 83336580
$ws.Range("L113").Value = 4478.7
$ws.Range("M113").Value = -83333326
$ws.Range("N113").Value = -10986.7
# Row 129
$ws.Range("H129").Value = 345602.97
$ws.Range("I129").Value = 333.33334
$ws.Range("J129").Value = 385441.78
$ws.Range("K129").Value = 1000.00002
$ws.Range("L129").Value = 1156325.34
$ws.Range("M129").Value = 3999.99998
$ws.Range("N129").Value = -1166325.34
# Row 132
$ws.Range("H132").Value = 3530.2856
$ws.Range("I132").Value = 3810.7273
$ws.Range("J132").Value = 2502
$ws.Range("K132").Value = 11432.1819
$ws.Range("L132").Value = 7506
$ws.Range("M132").Value = -8902.1819
$ws.Range("N132").Value = -12566
# Row 134
$ws.Range("H134").Value = 49900
$ws.Range("J134").Value = 49900
$ws.Range("L134").Value = 49900
$ws.Range("N134").Value = -60040
# Row 135
$ws.Range("H135").Value = 23819190
$ws.Range("I135").Value = 1230.7693
$ws.Range("K135").Value = 11076.9237
$ws.Range("M135").Value = -8541.923699999999
# Row 138
$ws.Range("H138").Value = 1409.1013
$ws.Range("I138").Value = 512.31915
$ws.Range("J138").Value = 2726.25
$ws.Range("K138").Value = 1536.95745
$ws.Range("L138").Value = 8178.75
$ws.Range("M138").Value = 3603.04255
$ws.Range("N138").Value = -18458.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1343.5938
$ws.Range("I2").Value = 1163.875
$ws.Range("K2").Value = 1163.875
$ws.Range("M2").Value = -1050.875
# Row 45
$ws.Range("H45").Value = 2461.35
$ws.Range("I45").Value = 2584.4167
$ws.Range("J45").Value = 2276.75
$ws.Range("K45").Value = 2584.4167
$ws.Range("L45").Value = 2276.75
$ws.Range("M45").Value = -2207.4167
$ws.Range("N45").Value = -3030.75
# Row 61
$ws.Range("H61").Value = 2911.1738
$ws.Range("I61").Value = 2417.4285
$ws.Range("J61").Value = 3679.2222
$ws.Range("K61").Value = 2417.4285
$ws.Range("L61").Value = 3679.2222
$ws.Range("M61").Value = -2205.4285
$ws.Range("N61").Value = -4103.2222
# Row 74
$ws.Range("H74").Value = 71429480
$ws.Range("I74").Value = 166667100
$ws.Range("K74").Value = 166667100
$ws.Range("M74").Value = -166666226
# Row 77
$ws.Range("H77").Value = 71429480
$ws.Range("I77").Value = 166667100
$ws.Range("K77").Value = 833335500
$ws.Range("M77").Value = -833331132
# Row 116
$ws.Range("H116").Value = 1343.5938
$ws.Range("I116").Value = 1163.875
$ws.Range("K116").Value = 1163.875
$ws.Range("M116").Value = 1130.125
# Row 136
$ws.Range("H136").Value = 2911.1738
$ws.Range("I136").Value = 2417.4285
$ws.Range("J136").Value = 3679.2222
$ws.Range("K136").Value = 7252.2855
$ws.Range("L136").Value = 11037.6666
$ws.Range("M136").Value = -4702.2855
$ws.Range("N136").Value = -16137.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1343.5938
$ws.Range("I3").Value = 1163.875
$ws.Range("K3").Value = 1163.875
$ws.Range("M3").Value = -1049.875
# Row 99
$ws.Range("H99").Value = 1453.4706
$ws.Range("I99").Value = 1190.909
$ws.Range("J99").Value = 1934.8334
$ws.Range("K99").Value = 1190.909
$ws.Range("L99").Value = 1934.8334
$ws.Range("M99").Value = 307.0909999999999
$ws.Range("N99").Value = -4930.8334

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 6000
$ws.Range("I33").Value = 6000
$ws.Range("K33").Value = 6000
$ws.Range("M33").Value = -5621
# Row 58
$ws.Range("H58").Value = 16258.091
$ws.Range("I58").Value = 1153.3043
$ws.Range("J58").Value = 50999.1
$ws.Range("K58").Value = 1153.3043
$ws.Range("L58").Value = 50999.1
$ws.Range("M58").Value = -950.3043
$ws.Range("N58").Value = -51405.1
# Row 99
$ws.Range("H99").Value = 21743500
$ws.Range("I99").Value = 3901.0908
$ws.Range("J99").Value = 41671464
$ws.Range("K99").Value = 3901.0908
$ws.Range("L99").Value = 41671464
$ws.Range("M99").Value = -2403.0908
$ws.Range("N99").Value = -41674460
# Row 126
$ws.Range("H126").Value = 21743500
$ws.Range("I126").Value = 3901.0908
$ws.Range("J126").Value = 41671464
$ws.Range("K126").Value = 11703.2724
$ws.Range("L126").Value = 125014392
$ws.Range("M126").Value = -9233.2724
$ws.Range("N126").Value = -125019332
# Row 134
$ws.Range("H134").Value = 1274.6046
$ws.Range("I134").Value = 956.9048
$ws.Range("J134").Value = 1577.8636
$ws.Range("K134").Value = 2870.7144
$ws.Range("L134").Value = 4733.5908
$ws.Range("M134").Value = -335.7143999999998
$ws.Range("N134").Value = -9803.5908
# Row 136
$ws.Range("H136").Value = 16258.091
$ws.Range("I136").Value = 1153.3043
$ws.Range("J136").Value = 50999.1
$ws.Range("K136").Value = 3459.9129
$ws.Range("L136").Value = 152997.3
$ws.Range("M136").Value = -909.9129000000003
$ws.Range("N136").Value = -158097.3

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 333965.84
$ws.Range("I129").Value = 531.6667
$ws.Range("J129").Value = 667400
$ws.Range("K129").Value = 1595.0001
$ws.Range("L129").Value = 2002200
$ws.Range("M129").Value = 3404.9999
$ws.Range("N129").Value = -2012200
# Row 131
$ws.Range("H131").Value = 754.66
$ws.Range("I131").Value = 539.8
$ws.Range("J131").Value = 765.96844
$ws.Range("K131").Value = 1619.4
$ws.Range("L131").Value = 2297.90532
$ws.Range("M131").Value = 3420.6
$ws.Range("N131").Value = -12377.90532
# Row 132
$ws.Range("H132").Value = 1501.2727
$ws.Range("I132").Value = 1404.8334
$ws.Range("J132").Value = 1617
$ws.Range("K132").Value = 12643.5006
$ws.Range("L132").Value = 14553
$ws.Range("M132").Value = -10113.5006
$ws.Range("N132").Value = -19613

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("M46").ClearContents()
$ws.Range("H46").Value = 25600
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 25600
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 25600
$ws.Range("N46").Value = -25912

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 5392.5386
$ws.Range("J61").Value = 7500.5
$ws.Range("L61").Value = 7500.5
$ws.Range("N61").Value = -7904.5
# Row 113
$ws.Range("H113").Value = 5392.5386
$ws.Range("J113").Value = 7500.5
$ws.Range("L113").Value = 7500.5
$ws.Range("N113").Value = -11840.5
